$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from A139 (bold, centered, bordered) onto the new index column cells A140:A149
$ws.Range("A139").Copy() | Out-Null
$ws.Range("A140:A149").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 140: Atalanta vs Milan
$ws.Cells.Item(140, 1).Value = 138
$ws.Cells.Item(140, 2).Value = "Atalanta"
$ws.Cells.Item(140, 3).Value = "Milan"
$ws.Cells.Item(140, 4).Value = 2
$ws.Cells.Item(140, 5).Value = 1
$ws.Cells.Item(140, 6).Value = 2.11
$ws.Cells.Item(140, 7).Value = 0.72
$ws.Cells.Item(140, 8).Value = 2.37
$ws.Cells.Item(140, 9).Value = 0.99
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 0.26
$ws.Cells.Item(140, 13).Value = 0.27
$ws.Cells.Item(140, 14).Value = 0.53
$ws.Cells.Item(140, 15).Value = 3

# Row 141: Inter vs Parma
$ws.Cells.Item(141, 1).Value = 139
$ws.Cells.Item(141, 2).Value = "Inter"
$ws.Cells.Item(141, 3).Value = "Parma"
$ws.Cells.Item(141, 4).Value = 3
$ws.Cells.Item(141, 5).Value = 1
$ws.Cells.Item(141, 6).Value = 1.77
$ws.Cells.Item(141, 7).Value = 0.38
$ws.Cells.Item(141, 8).Value = 2.12
$ws.Cells.Item(141, 9).Value = 0.43
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 0.35
$ws.Cells.Item(141, 13).Value = 0.05
$ws.Cells.Item(141, 14).Value = 0.4
$ws.Cells.Item(141, 15).Value = 4

# Row 142: Genoa vs Torino
$ws.Cells.Item(142, 1).Value = 140
$ws.Cells.Item(142, 2).Value = "Genoa"
$ws.Cells.Item(142, 3).Value = "Torino"
$ws.Cells.Item(142, 4).Value = 0
$ws.Cells.Item(142, 5).Value = 0
$ws.Cells.Item(142, 6).Value = 0.49
$ws.Cells.Item(142, 7).Value = 0.59
$ws.Cells.Item(142, 8).Value = 0.47
$ws.Cells.Item(142, 9).Value = 0.58
$ws.Cells.Item(142, 10).Value = 0
$ws.Cells.Item(142, 11).Value = 0
$ws.Cells.Item(142, 12).Value = 0.02
$ws.Cells.Item(142, 13).Value = 0.01
$ws.Cells.Item(142, 14).Value = 0.03
$ws.Cells.Item(142, 15).Value = 0

# Row 143: Juventus vs Bologna
$ws.Cells.Item(143, 1).Value = 141
$ws.Cells.Item(143, 2).Value = "Juventus"
$ws.Cells.Item(143, 3).Value = "Bologna"
$ws.Cells.Item(143, 4).Value = 2
$ws.Cells.Item(143, 5).Value = 2
$ws.Cells.Item(143, 6).Value = 1.02
$ws.Cells.Item(143, 7).Value = 0.93
$ws.Cells.Item(143, 8).Value = 0.98
$ws.Cells.Item(143, 9).Value = 0.58
$ws.Cells.Item(143, 10).Value = 0
$ws.Cells.Item(143, 11).Value = 0
$ws.Cells.Item(143, 12).Value = 0.04
$ws.Cells.Item(143, 13).Value = 0.35
$ws.Cells.Item(143, 14).Value = 0.38
$ws.Cells.Item(143, 15).Value = 4

# Row 144: Roma vs Lecce
$ws.Cells.Item(144, 1).Value = 142
$ws.Cells.Item(144, 2).Value = "Roma"
$ws.Cells.Item(144, 3).Value = "Lecce"
$ws.Cells.Item(144, 4).Value = 4
$ws.Cells.Item(144, 5).Value = 1
$ws.Cells.Item(144, 6).Value = 1.9
$ws.Cells.Item(144, 7).Value = 0.88
$ws.Cells.Item(144, 8).Value = 2.51
$ws.Cells.Item(144, 9).Value = 0.9399999999999999
$ws.Cells.Item(144, 10).Value = 0
$ws.Cells.Item(144, 11).Value = 1
$ws.Cells.Item(144, 12).Value = 0.61
$ws.Cells.Item(144, 13).Value = 0.06
$ws.Cells.Item(144, 14).Value = 0.67
$ws.Cells.Item(144, 15).Value = 4

# Row 145: Fiorentina vs Cagliari
$ws.Cells.Item(145, 1).Value = 143
$ws.Cells.Item(145, 2).Value = "Fiorentina"
$ws.Cells.Item(145, 3).Value = "Cagliari"
$ws.Cells.Item(145, 4).Value = 1
$ws.Cells.Item(145, 5).Value = 0
$ws.Cells.Item(145, 6).Value = 0.32
$ws.Cells.Item(145, 7).Value = 0.84
$ws.Cells.Item(145, 8).Value = 0.42
$ws.Cells.Item(145, 9).Value = 1.1
$ws.Cells.Item(145, 10).Value = 0
$ws.Cells.Item(145, 11).Value = 1
$ws.Cells.Item(145, 12).Value = 0.1
$ws.Cells.Item(145, 13).Value = 0.26
$ws.Cells.Item(145, 14).Value = 0.36
$ws.Cells.Item(145, 15).Value = 2

# Row 146: Hellas Verona vs Empoli
$ws.Cells.Item(146, 1).Value = 144
$ws.Cells.Item(146, 2).Value = "Hellas Verona"
$ws.Cells.Item(146, 3).Value = "Empoli"
$ws.Cells.Item(146, 4).Value = 1
$ws.Cells.Item(146, 5).Value = 4
$ws.Cells.Item(146, 6).Value = 0.9399999999999999
$ws.Cells.Item(146, 7).Value = 1.31
$ws.Cells.Item(146, 8).Value = 0.95
$ws.Cells.Item(146, 9).Value = 1.34
$ws.Cells.Item(146, 10).Value = 0
$ws.Cells.Item(146, 11).Value = 0
$ws.Cells.Item(146, 12).Value = 0.01
$ws.Cells.Item(146, 13).Value = 0.03
$ws.Cells.Item(146, 14).Value = 0.04
$ws.Cells.Item(146, 15).Value = 5

# Row 147: Napoli vs Lazio
$ws.Cells.Item(147, 1).Value = 145
$ws.Cells.Item(147, 2).Value = "Napoli"
$ws.Cells.Item(147, 3).Value = "Lazio"
$ws.Cells.Item(147, 4).Value = 0
$ws.Cells.Item(147, 5).Value = 1
$ws.Cells.Item(147, 6).Value = 0.5600000000000001
$ws.Cells.Item(147, 7).Value = 0.29
$ws.Cells.Item(147, 8).Value = 0.63
$ws.Cells.Item(147, 9).Value = 0.37
$ws.Cells.Item(147, 10).Value = 0
$ws.Cells.Item(147, 11).Value = 0
$ws.Cells.Item(147, 12).Value = 0.07000000000000001
$ws.Cells.Item(147, 13).Value = 0.08
$ws.Cells.Item(147, 14).Value = 0.16
$ws.Cells.Item(147, 15).Value = 1

# Row 148: Venezia vs Como
$ws.Cells.Item(148, 1).Value = 146
$ws.Cells.Item(148, 2).Value = "Venezia"
$ws.Cells.Item(148, 3).Value = "Como"
$ws.Cells.Item(148, 4).Value = 2
$ws.Cells.Item(148, 5).Value = 2
$ws.Cells.Item(148, 6).Value = 0.33
$ws.Cells.Item(148, 7).Value = 1.05
$ws.Cells.Item(148, 8).Value = 0.43
$ws.Cells.Item(148, 9).Value = 0.85
$ws.Cells.Item(148, 10).Value = 0
$ws.Cells.Item(148, 11).Value = 0
$ws.Cells.Item(148, 12).Value = 0.1
$ws.Cells.Item(148, 13).Value = 0.2
$ws.Cells.Item(148, 14).Value = 0.3
$ws.Cells.Item(148, 15).Value = 4

# Row 149: Monza vs Udinese
$ws.Cells.Item(149, 1).Value = 147
$ws.Cells.Item(149, 2).Value = "Monza"
$ws.Cells.Item(149, 3).Value = "Udinese"
$ws.Cells.Item(149, 4).Value = 1
$ws.Cells.Item(149, 5).Value = 2
$ws.Cells.Item(149, 6).Value = 1.32
$ws.Cells.Item(149, 7).Value = 1.22
$ws.Cells.Item(149, 8).Value = 2.06
$ws.Cells.Item(149, 9).Value = 1.37
$ws.Cells.Item(149, 10).Value = 0
$ws.Cells.Item(149, 11).Value = 0
$ws.Cells.Item(149, 12).Value = 0.74
$ws.Cells.Item(149, 13).Value = 0.15
$ws.Cells.Item(149, 14).Value = 0.88
$ws.Cells.Item(149, 15).Value = 3

Write-Host "Added rows 140-149 (round 15 matches)"
